$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.528.66'
$ws.Range("E2").Value = '  +3.09%  '
$ws.Range("D3").Value = '3.708.17'
$ws.Range("E3").Value = '  +7.97%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '584.41'
$ws.Range("E5").Value = '  +0.63%  '
$ws.Range("D6").Value = '179.59'
$ws.Range("E6").Value = '  +1.20%  '
$ws.Range("D7").Value = '3.697.90'
$ws.Range("E7").Value = '  +7.86%  '
$ws.Range("D8").Value = '0.615'
$ws.Range("E8").Value = '  +3.98%  '
$ws.Range("E9").Value = '  +0.07%  '
$ws.Range("D10").Value = '0.201'
$ws.Range("E10").Value = '  +1.53%  '
$ws.Range("D11").Value = '0.611'
$ws.Range("E11").Value = '  +4.34%  '
$ws.Range("D12").Value = '49.41'
$ws.Range("E12").Value = '  +1.21%  '
$ws.Range("E13").Value = '  +2.53%  '
$ws.Range("D14").Value = '4.305.90'
$ws.Range("E14").Value = '  +8.56%  '
$ws.Range("D15").Value = '679.07'
$ws.Range("E15").Value = '  -2.92%  '
$ws.Range("D16").Value = '9.06'
$ws.Range("E16").Value = '  +4.78%  '
$ws.Range("D17").Value = '3.715.45'
$ws.Range("E17").Value = '  +8.18%  '
$ws.Range("D18").Value = '71.589.73'
$ws.Range("E18").Value = '  +3.14%  '
$ws.Range("D19").Value = '0.122'
$ws.Range("E19").Value = '  +1.18%  '
$ws.Range("D20").Value = '18.04'
$ws.Range("E20").Value = '  +1.71%  '
$ws.Range("E21").Value = '  +1.94%  '
$ws.Range("E22").Value = '  +18.32%  '
$ws.Range("D23").Value = '0.944'
$ws.Range("E23").Value = '  +4.94%  '
$ws.Range("D24").Value = '17.53'
$ws.Range("E24").Value = '  +3.23%  '
$ws.Range("D25").Value = '102.55'
$ws.Range("E25").Value = '  +1.36%  '
$ws.Range("E26").Value = '  +2.87%  '
$ws.Range("D27").Value = '2.84'
$ws.Range("E27").Value = '  +5.98%  '
$ws.Range("E28").Value = '  +8.92%  '
$ws.Range("D29").Value = '35.52'
$ws.Range("E29").Value = '  +5.63%  '
$ws.Range("D30").Value = '9.21'
$ws.Range("E30").Value = '  +5.04%  '
$ws.Range("D31").Value = '7.36'
$ws.Range("E31").Value = '  +5.32%  '
$ws.Range("D32").Value = '4.13'
$ws.Range("E32").Value = '  +9.69%  '
$ws.Range("D33").Value = '590.97'
$ws.Range("E33").Value = '  +4.22%  '
$ws.Range("D34").Value = '11.22'
$ws.Range("E34").Value = '  +1.66%  '
$ws.Range("E35").Value = '  +3.50%  '
$ws.Range("E36").Value = '  +1.90%  '
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").Value = '3.683.78'
$ws.Range("E38").Value = '  +2.02%  '
$ws.Range("E39").Value = '  +5.02%  '
$ws.Range("D40").Value = '0.0₃0779'
$ws.Range("E40").Value = '  +6.51%  '
$ws.Range("D41").Value = '35.71'
$ws.Range("E41").Value = '  +2.05%  '
$ws.Range("D42").Value = '3.45'
$ws.Range("E42").Value = '  +4.90%  '
$ws.Range("D43").Value = '2.79'
$ws.Range("E43").Value = '  +4.31%  '
$ws.Range("D44").Value = "'0.0460"
$ws.Range("E44").Value = '  +9.68%  '
$ws.Range("E45").Value = '  +4.84%  '
$ws.Range("D46").Value = '2.88'
$ws.Range("E46").Value = '  +8.38%  '
$ws.Range("D47").Value = '3.39'
$ws.Range("E47").Value = '  +1.19%  '
$ws.Range("E48").Value = '  +4.11%  '
$ws.Range("D49").Value = '1.45'
$ws.Range("E49").Value = '  -1.31%  '
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("D51").Value = '135.64'
$ws.Range("E51").Value = '  +3.13%  '
